# Apply the commit: add a new "LikelihoodCalcEx" worksheet in front of the
# existing sheets, containing a small worked example of the gradient /
# likelihood calculation referenced in the commit message
# ("added gradient calculation to manual").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new sheet as the very first tab.
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "LikelihoodCalcEx"

# ---------------------------------------------------------------------
# 2. Column widths (matches the style widths used on the sheet).
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 13
$ws.Range("B1:D1").EntireColumn.ColumnWidth = 10.83203125
$ws.Columns.Item(5).ColumnWidth = 9.33203125
$ws.Columns.Item(6).ColumnWidth = 5.33203125
$ws.Columns.Item(7).ColumnWidth = 1.1640625
$ws.Columns.Item(8).ColumnWidth = 7
$ws.Range("I1:J1").EntireColumn.ColumnWidth = 9.6640625
$ws.Columns.Item(11).ColumnWidth = 10.1640625
$ws.Columns.Item(12).ColumnWidth = 13

# ---------------------------------------------------------------------
# 3. Header / parameter block (rows 1-5).
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "sigma:"
$ws.Range("A3").Value = 100

$ws.Range("B2").Value = "L1 penalty"
$ws.Range("B3").Value = "mu:"
$ws.Range("B4").Value = "L2 penalty"
$ws.Range("B5").Value = "weights -->"

$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0

$ws.Range("C5").Value = 1.1
$ws.Range("D5").Value = 2.3
$ws.Range("E5").Value = 4.8

$ws.Range("C2").Formula = "=(C5-C3)*SQRT(2)/`$A`$3^2"
$ws.Range("D2").Formula = "=(D5-D3)*SQRT(2)/`$A`$3^2"
$ws.Range("E2").Formula = "=(E5-E3)*SQRT(2)/`$A`$3^2"

$ws.Range("C4").Formula = "=(C5-C3)^2/2/`$A`$3"
$ws.Range("D4").Formula = "=(D5-D3)^2/2/`$A`$3"
$ws.Range("E4").Formula = "=(E5-E3)^2/2/`$A`$3"

$ws.Range("I2").Formula = "=LN(1)"

# ---------------------------------------------------------------------
# 4. Table header row (row 6).
# ---------------------------------------------------------------------
$ws.Range("A6").Value = "input"
$ws.Range("B6").Value = "candidate"
$ws.Range("C6").Value = "C1"
$ws.Range("D6").Value = "C2"
$ws.Range("E6").Value = "C3"
$ws.Range("F6").Value = "H"
$ws.Range("G6").Value = "eH"
$ws.Range("H6").Value = "cand.prob"
$ws.Range("I6").Value = "Likelihood:"
$ws.Range("J6").Formula = "=SUMPRODUCT(J9)"
$ws.Range("K6").Value = "Observed"

# ---------------------------------------------------------------------
# 5. Data rows (7-11).
# ---------------------------------------------------------------------
$ws.Range("A7").Value = "input_1"
$ws.Range("B7").Value = "cand_1"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 0
$ws.Range("K7").Value = 415

$ws.Range("A8").Value = "input_1"
$ws.Range("B8").Value = "cand_1"
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 1
$ws.Range("K8").Value = 12

$ws.Range("A9").Value = "input_1"
$ws.Range("B9").Value = "cand_1"
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 0
$ws.Range("K9").Value = 891

$ws.Range("A10").Value = "input_2"
$ws.Range("B10").Value = "cand_2"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 1
$ws.Range("K10").Value = 1

$ws.Range("A11").Value = "input_2"
$ws.Range("B11").Value = "cand_2"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("K11").Value = 345

# F/G/H/I formulas for rows 7-11 (F & H differ per row only by sumif range,
# G & I are uniform shared-style formulas).
$ws.Range("F7").Formula = "=-SUMPRODUCT(C`$5:E`$5,C7:E7)"
$ws.Range("F8").Formula = "=-SUMPRODUCT(C`$5:E`$5,C8:E8)"
$ws.Range("F9").Formula = "=-SUMPRODUCT(C`$5:E`$5,C9:E9)"
$ws.Range("F10").Formula = "=-SUMPRODUCT(C`$5:E`$5,C10:E10)"
$ws.Range("F11").Formula = "=-SUMPRODUCT(C`$5:E`$5,C11:E11)"

$ws.Range("G7:G11").Formula = "=EXP(F7)"

$ws.Range("H7").Formula = "=G7/SUMIF(A:A,A7,G:G)"
$ws.Range("H8:H11").Formula = "=G8/SUMIF(A:A,A8,G:G)"

$ws.Range("I7:I11").Formula = "=LN(H7)"

# ---------------------------------------------------------------------
# 6. Formatting.
# ---------------------------------------------------------------------
# Whole used range: apply font/fill "pass-through" so every cell carries
# applyFont/applyFill like the authored sheet.
$used = $ws.Range("A1:K11")
$used.Font.Name = "Calibri"
$used.Font.Size = 12

# Header row 6 styling (bold, centered, light border row).
$hdr = $ws.Range("A6:K6")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108

# Numeric columns centered.
$ws.Range("C1:E11").HorizontalAlignment = -4108
$ws.Range("F1:F11").HorizontalAlignment = -4108
$ws.Range("G1:G11").HorizontalAlignment = -4108
$ws.Range("H1:H11").HorizontalAlignment = -4108
$ws.Range("I1:I11").HorizontalAlignment = -4108
$ws.Range("K1:K11").HorizontalAlignment = -4108

# C2:E2 and C4:E4 use 2-decimal number format (numFmtId 2).
$ws.Range("C2:E2").NumberFormat = "0.00"
$ws.Range("C4:E4").NumberFormat = "0.00"

# Bottom border under row 5 (weights row) to set off the header below it.
$ws.Range("A5:K5").Borders.Item(9).LineStyle = 1

# Row 10 is visually distinguished (custom row format) in the source sheet.
$ws.Rows.Item(10).Font.Name = "Calibri"
$ws.Rows.Item(10).Font.Size = 12

$ws.Range("E6").Select()
